# K-Fold cross validation testing
# Adds a "10-Fold" results sheet and appends Training/Testing time rows
# to the "initial" sheet, matching a 10-fold cross-validation run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "10-Fold" sheet after "scaled" (the last sheet) and
#    populate the K-Fold cross-validation summary. Done first so that
#    shared strings introduced here ("Training Time") land ahead of the
#    ones introduced on the "initial" sheet below, matching the order
#    they were authored in.
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$wsFold = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$wsFold.Name = "10-Fold"

$wsFold.Range("B1").Value = "Random Forect"
$wsFold.Range("B1").Font.Bold = $true
$wsFold.Range("C1").Value = "SVM"
$wsFold.Range("C1").Font.Bold = $true
$wsFold.Range("D1").Value = "Logistic Regression"
$wsFold.Range("D1").Font.Bold = $true
$wsFold.Range("E1").Value = "GaussianNB"
$wsFold.Range("E1").Font.Bold = $true
$wsFold.Range("F1").Value = "KNeighborsClassifier"
$wsFold.Range("F1").Font.Bold = $true
$wsFold.Range("G1").Value = "Decision Tree"
$wsFold.Range("G1").Font.Bold = $true
$wsFold.Range("H1").Value = "AdaBoost"
$wsFold.Range("H1").Font.Bold = $true
$wsFold.Range("I1").Value = "GradientBoosting"
$wsFold.Range("I1").Font.Bold = $true
$wsFold.Range("J1").Value = "Bagging"
$wsFold.Range("J1").Font.Bold = $true

$wsFold.Range("A2").Value = "Accuracy"
$wsFold.Range("A2").Font.Bold = $true
$wsFold.Range("B2").Value = 0.93
$wsFold.Range("D2").Value = 0.703
$wsFold.Range("E2").Value = 0.696
$wsFold.Range("F2").Value = 0.811
$wsFold.Range("G2").Value = 0.708
$wsFold.Range("H2").Value = 0.705
$wsFold.Range("I2").Value = 0.729
$wsFold.Range("J2").Value = 0.703

$wsFold.Range("A4").Value = "Training Time"
$wsFold.Range("A4").Font.Bold = $true
$wsFold.Range("B4").Value = 317.1
$wsFold.Range("D4").Value = 10.8
$wsFold.Range("E4").Value = 1.74
$wsFold.Range("F4").Value = 33.1
$wsFold.Range("G4").Value = 10.4
$wsFold.Range("H4").Value = 194.53
$wsFold.Range("I4").Value = 1042.7
$wsFold.Range("J4").Value = 7134

$wsFold.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 2) "initial" sheet (sheet1): append Training Time(s) / Testing Time(s)
#    rows (16 & 17) below the existing metrics table.
# ---------------------------------------------------------------------
$wsInitial = $wb.Worksheets.Item("initial")

$wsInitial.Range("A16").Value = "Training Time(s)"
$wsInitial.Range("A16").Font.Bold = $true
$wsInitial.Range("B16").Value = 16.92
$wsInitial.Range("C16").Value = 346.7
$wsInitial.Range("D16").Value = 0.067
$wsInitial.Range("E16").Value = 0.07
$wsInitial.Range("F16").Value = 0.03
$wsInitial.Range("G16").Value = 0.62
$wsInitial.Range("H16").Value = 10.92
$wsInitial.Range("I16").Value = 52.74
$wsInitial.Range("J16").Value = 265.75

$wsInitial.Range("A17").Value = "Testing Time(s)"
$wsInitial.Range("A17").Font.Bold = $true
$wsInitial.Range("B17").Value = 0.44
$wsInitial.Range("C17").Value = 9.77
$wsInitial.Range("D17").Value = 0.02
$wsInitial.Range("E17").Value = 0.03
$wsInitial.Range("F17").Value = 7.3
$wsInitial.Range("G17").Value = 0.02
$wsInitial.Range("H17").Value = 0.26
$wsInitial.Range("I17").Value = 0.05
$wsInitial.Range("J17").Value = 187.96

# Restore the selection that was active on this sheet before the new
# 10-Fold tab became active.
$wsInitial.Range("B12").Select()

# ---------------------------------------------------------------------
# 3) Leave the new "10-Fold" sheet active/selected, as it was the last
#    sheet the author worked on.
# ---------------------------------------------------------------------
$wsFold.Activate()
$wsFold.Range("J4").Select()
